$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab from "Sheet1" to "Sheet 1"
$ws.Name = "Sheet 1"

# Insert a new column before column B to hold the "Reef" grouping (C/U/Z),
# pushing Lat_Y/Long_X from B:C to C:D.
$ws.Columns("B:B").Insert()

# Header for the new column
$ws.Range("B1").Value = "Reef"

# Reef codes per site (IC-C1/IC-C2 -> C, IC-U1/IC-U2 -> U, IC-Z1/IC-Z2 -> Z)
$ws.Range("B2").Value = "C"
$ws.Range("B3").Value = "C"
$ws.Range("B4").Value = "U"
$ws.Range("B5").Value = "U"
$ws.Range("B6").Value = "Z"
$ws.Range("B7").Value = "Z"

# Update the selection to match the saved view (B5 active cell)
$ws.Range("B5").Select()
